$wb = $excel.ActiveWorkbook

# --- "Games" sheet: append the completed game (2024-01-15 vs MEM) as new row 41 ---
$games = $wb.Worksheets.Item("Games")

$games.Cells.Item(41, 1).Value = 40
$games.Cells.Item(41, 2).Value = 45306
$games.Cells.Item(41, 2).NumberFormat = "YYYY-MM-DD"
$games.Cells.Item(41, 3).Value = -2
$games.Cells.Item(41, 4).Value = 107
$games.Cells.Item(41, 5).Value = 102
$games.Cells.Item(41, 6).Value = 0.544
$games.Cells.Item(41, 7).Value = 16.8
$games.Cells.Item(41, 8).Value = 18.4
$games.Cells.Item(41, 9).Value = 0.1
$games.Cells.Item(41, 10).Value = 104.9
$games.Cells.Item(41, 11).Value = "MEM"
$games.Cells.Item(41, 12).Value = 116
$games.Cells.Item(41, 13).Value = 0.488
$games.Cells.Item(41, 14).Value = 9.6
$games.Cells.Item(41, 15).Value = 22.4
$games.Cells.Item(41, 16).Value = 0.372
$games.Cells.Item(41, 17).Value = 113.7
$games.Cells.Item(41, 18).Value = 0
$games.Cells.Item(41, 19).Value = 0

# --- "Next" sheet: that fixture has now been played, so drop it from the schedule ---
$next = $wb.Worksheets.Item("Next")
$next.Rows.Item(2).Delete()
